$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "WithTable" - table MyTable at A1:E5 -> A1:F5. A new "Boolean" column
# is inserted as the 4th column (between "DateTime" and "ARCtrl Column"), so
# the two "ARCtrl Column" columns shift one column to the right.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("WithTable")
$tbl1 = $ws1.ListObjects.Item(1)

# Insert a blank column before D, shifting the two "ARCtrl Column" columns right
$ws1.Range("D1:D5").EntireColumn.Insert(-4161)

# Header
$ws1.Range("D1").Value = "Boolean"

# Data: copy the DateTime column's style (same number-format style as the
# date cells) down into the new column, then overwrite with boolean literals
# (keeps the style but switches the stored type to boolean true/false).
$ws1.Range("C2:C5").Copy($ws1.Range("D2:D5"))
$ws1.Range("D2").Value = $true
$ws1.Range("D3").Value = $false
$ws1.Range("D4").Value = $true
$ws1.Range("D5").Value = $false

# Grow the table to include the new column plus the now-shifted last column
$tbl1.Resize($ws1.Range("A1:F5"))

# Re-stamp every header cell so the table's column-name metadata re-syncs
# with the (already-correct) header text for every column, not just the
# newest one.
foreach ($colLetter in @("A","B","C","D","E","F")) {
    $cell = $ws1.Range($colLetter + "1")
    $v = $cell.Value2
    $cell.Value = $v
}

$ws1.Range("F5").Select()

# ---------------------------------------------------------------------------
# Sheet "Tableless" - same underlying data grid, no table object attached.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Tableless")

$ws2.Range("D1:D5").EntireColumn.Insert(-4161)

$ws2.Range("D1").Value = "Boolean"

$ws2.Range("C2:C5").Copy($ws2.Range("D2:D5"))
$ws2.Range("D2").Value = $true
$ws2.Range("D3").Value = $false
$ws2.Range("D4").Value = $true
$ws2.Range("D5").Value = $false

$ws2.Range("D6").Select()

# ---------------------------------------------------------------------------
# Sheet "WithTable_Duplicate" - table MyOtherTable at B4:F8 -> B4:G8. Data
# starts one column to the right (column B) and four rows down (row 4), so
# the new "Boolean" column goes in column E.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("WithTable_Duplicate")
$tbl2 = $ws3.ListObjects.Item(1)

$ws3.Range("E4:E8").EntireColumn.Insert(-4161)

$ws3.Range("E4").Value = "Boolean"

$ws3.Range("D5:D8").Copy($ws3.Range("E5:E8"))
$ws3.Range("E5").Value = $true
$ws3.Range("E6").Value = $false
$ws3.Range("E7").Value = $true
$ws3.Range("E8").Value = $false

$tbl2.Resize($ws3.Range("B4:G8"))

foreach ($colLetter in @("B","C","D","E","F","G")) {
    $cell = $ws3.Range($colLetter + "4")
    $v = $cell.Value2
    $cell.Value = $v
}

$ws3.Range("F7").Select()
